$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'67.157.27"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  -4.18%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'3.241.56"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  -7.57%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('D4').Value = "'0.999"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = "'  -0.07%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'596.34"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  -1.42%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'153.51"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  -11.07%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('E7').Value = "'  -0.06%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = "'3.233.92"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'  -7.64%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('E9').Value = "'  -10.04%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'0.174"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  -10.18%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'6.66"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'  -8.16%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'0.504"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'  -14.11%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'39.43"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  -14.42%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'0.0000248"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  -10.01%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'3.757.66"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  -7.75%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'67.152.91"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  -4.20%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'3.242.25"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  -7.55%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('E18').Value = "'  -4.61%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('E19').Value = "'  -13.59%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = "'534.37"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'  -13.04%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('E21').Value = "'  -13.62%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('E22').Value = "'  -12.78%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'7.93"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  -13.49%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'13.95"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  -10.16%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'86.22"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'  -12.58%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('E26').Value = "'  -0.04%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('E27').Value = "'  -14.19%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('E28').Value = "'  -13.63%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = "'8.23"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'  -9.00%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = "'29.53"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'  -12.54%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('E31').Value = "'  -11.76%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = "'1.16"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'  -8.90%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('B33').Value = "'Bittensor"
$ws.Range('B33').Style = 'Normal'
$ws.Range('C33').Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range('C33').Style = 'Normal'
$ws.Range('D33').Value = "'540.66"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'  -14.26%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('B34').Value = "'Filecoin"
$ws.Range('B34').Style = 'Normal'
$ws.Range('C34').Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range('C34').Style = 'Normal'
$ws.Range('D34').Value = "'6.61"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'  -17.92%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('E35').Value = "'  -14.93%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('E36').Value = "'  +0.18%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = "'53.27"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'  -6.13%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('E38').Value = "'  -12.27%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = "'0.0429"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'  -11.07%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = "'9.38"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'  -12.79%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('D42').Value = "'2.944.22"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'  -12.40%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('E43').Value = "'  -22.83%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('E44').Value = "'  -13.54%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('E45').Value = "'  -18.78%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('E46').Value = "'  -16.00%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'26.74"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  -16.32%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('E49').Value = "'  -0.16%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('E50').Value = "'  -11.61%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = "'123.11"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  -7.63%  "
$ws.Range('E51').Style = 'Normal'
